$d = $word.ActiveDocument

# The trailing "_GoBack" bookmark currently sits right after "fucntion" in the
# last paragraph. We will rebuild that paragraph (and add several new ones)
# via InsertXML, so remove the existing bookmark first -- we re-create it in
# its new resting place (end of the very last paragraph) as part of the XML
# fragment below.
$bm = $d.Bookmarks.Item("_GoBack")
if ($bm.Start -ne $null) {
    $bm.Delete()
}

# Locate the paragraph to rewrite by its current (typo-laden) text.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Created a DTO to populate a agrigate fucntion")
if (-not $found) {
    throw "Could not find the 'Created a DTO...' paragraph text"
}

# Re-seat the found bounds into a plain Range (a Range that has gone through
# Find.Execute behaves like an insertion point for InsertXML instead of a
# replacement span), then replace that whole span with the corrected
# paragraph plus the new October 5th content.
$target = $d.Range($searchRange.Start, $searchRange.End)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">Created a DTO to populate </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>a</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>aggregate</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>function</w:t></w:r>
</w:p>
<w:p/>
<w:p>
<w:pPr>
<w:rPr>
<w:vertAlign w:val="superscript"/>
</w:rPr>
</w:pPr>
<w:r><w:t>October 5</w:t></w:r>
<w:r>
<w:rPr>
<w:vertAlign w:val="superscript"/>
</w:rPr>
<w:t>th</w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:rPr>
<w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
<w:color w:val="808080"/>
<w:sz w:val="19"/>
<w:szCs w:val="19"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
<w:color w:val="808080"/>
<w:sz w:val="19"/>
<w:szCs w:val="19"/>
<w:highlight w:val="black"/>
</w:rPr>
<w:t>&lt;</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
<w:color w:val="569CD6"/>
<w:sz w:val="19"/>
<w:szCs w:val="19"/>
<w:highlight w:val="black"/>
</w:rPr>
<w:t>div</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
<w:color w:val="DCDCDC"/>
<w:sz w:val="19"/>
<w:szCs w:val="19"/>
<w:highlight w:val="black"/>
</w:rPr>
<w:t xml:space="preserve"> </w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
<w:color w:val="9CDCFE"/>
<w:sz w:val="19"/>
<w:szCs w:val="19"/>
<w:highlight w:val="black"/>
</w:rPr>
<w:t>class</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
<w:color w:val="B4B4B4"/>
<w:sz w:val="19"/>
<w:szCs w:val="19"/>
<w:highlight w:val="black"/>
</w:rPr>
<w:t>=</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
<w:color w:val="C8C8C8"/>
<w:sz w:val="19"/>
<w:szCs w:val="19"/>
<w:highlight w:val="black"/>
</w:rPr>
<w:t>"row col-md-12"</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
<w:color w:val="808080"/>
<w:sz w:val="19"/>
<w:szCs w:val="19"/>
<w:highlight w:val="black"/>
</w:rPr>
<w:t>&gt;</w:t>
</w:r>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">Using </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>ready made</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> CSS classes from Bootstrap</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.InsertXML($xml)
